$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 26 data. Force column A to be treated as plain text so the
# date-like string "2025-09-10" is not auto-converted into a date serial,
# then clear the temporary formatting so no extra style is left on the cell.
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "2025-09-10"
$ws.Cells.Item(26, 1).ClearFormats()

$ws.Cells.Item(26, 2).Value = 57.56999969482422
$ws.Cells.Item(26, 3).Value = 709.0999755859375
$ws.Cells.Item(26, 4).Value = 324.3999938964844
